$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.744.59'
$ws.Cells.Item(2, 5).Value = '  -3.92%  '
$ws.Cells.Item(3, 4).Value = '3.338.42'
$ws.Cells.Item(3, 5).Value = '  -4.77%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '182.27'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  -9.07%  '
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '533.53'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -3.21%  '
$origStyle = $ws.Cells.Item(7, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.609'
$ws.Cells.Item(7, 4).Style = $origStyle
$ws.Cells.Item(7, 5).Value = '  +0.38%  '
$ws.Cells.Item(8, 4).Value = '3.332.70'
$ws.Cells.Item(8, 5).Value = '  -4.81%  '
$ws.Cells.Item(9, 5).Value = '  +0.08%  '
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.619'
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  -5.44%  '
$origStyle = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '59.14'
$ws.Cells.Item(11, 4).Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  -7.12%  '
$origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.135'
$ws.Cells.Item(12, 4).Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  -5.29%  '
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000262'
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  -2.73%  '
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '9.20'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  -6.47%  '
$ws.Cells.Item(15, 4).Value = '3.872.54'
$ws.Cells.Item(15, 5).Value = '  -4.72%  '
$ws.Cells.Item(16, 4).Value = '3.343.40'
$ws.Cells.Item(16, 5).Value = '  -4.64%  '
$ws.Cells.Item(17, 5).Value = '  -4.35%  '
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '64.764.46'
$ws.Cells.Item(18, 5).Value = '  -3.52%  '
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.69'
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  -3.57%  '
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '11.26'
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  -4.38%  '
$origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.972'
$ws.Cells.Item(21, 4).Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  -4.98%  '
$origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '378.03'
$ws.Cells.Item(22, 4).Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  -3.40%  '
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.84'
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -3.99%  '
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.31'
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -7.13%  '
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '81.34'
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  -1.46%  '
$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.96'
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  +2.20%  '
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '6.09'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  -0.97%  '
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.70'
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  -3.62%  '
$origStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '11.57'
$ws.Cells.Item(29, 4).Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  -5.34%  '
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.48'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -3.84%  '
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '29.23'
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -5.58%  '
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '659.40'
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  -2.63%  '
$origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.77'
$ws.Cells.Item(33, 4).Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  -3.16%  '
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '11.38'
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  -3.12%  '
$origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.107'
$ws.Cells.Item(35, 4).Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  -2.92%  '
$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '59.77'
$ws.Cells.Item(36, 4).Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  -6.39%  '
$ws.Cells.Item(37, 2).Value = 'Dai'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.999'
$ws.Cells.Item(37, 4).Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  -0.08%  '
$ws.Cells.Item(38, 2).Value = 'TheGraph'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$origStyle = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.397'
$ws.Cells.Item(38, 4).Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  -0.33%  '
$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '37.19'
$ws.Cells.Item(39, 4).Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  -3.70%  '
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$origStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 4).Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  +0.14%  '
$ws.Cells.Item(41, 2).Value = 'PEPE'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(41, 4).Value = '0.0₃0715'
$ws.Cells.Item(41, 5).Value = '  +6.33%  '
$ws.Cells.Item(42, 5).Value = '  -3.03%  '
$ws.Cells.Item(43, 4).Value = '2.946.87'
$ws.Cells.Item(43, 5).Value = '  -4.05%  '
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.53'
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  +0.35%  '
$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.74'
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -7.91%  '
$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0403'
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  +1.46%  '
$ws.Cells.Item(47, 5).Value = '  -3.67%  '
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.09'
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  +6.57%  '
$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.82'
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  +7.92%  '
$ws.Cells.Item(50, 5).Value = '  +0.29%  '
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.54'
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  -4.88%  '
